$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 111, pushing existing rows 111-117 down to 112-118
$ws.Rows(111).Insert()

# Populate the new row 111 with the latest weekly entry
$ws.Range("A111").Value = 6
$ws.Range("B111").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C111").Value = "Metropolitana"
$ws.Range("D111").Value = 44706
$ws.Range("E111").Value = 13
$ws.Range("F111").Value = "Fruta"
$ws.Range("G111").Value = 100104
$ws.Range("H111").Value = "Frutos de pepita"
$ws.Range("I111").Value = 100104003
$ws.Range("J111").Value = "Membrillo"
$ws.Range("K111").Value = "Champion"
$ws.Range("L111").Value = "Primera"
$ws.Range("M111").Value = 15
$ws.Range("N111").Value = 230000
$ws.Range("O111").Value = 230000
$ws.Range("P111").Value = 230000
$ws.Range("Q111").Value = "$/bins (450 kilos)"
$ws.Range("R111").Value = "Región Metropolitana"
$ws.Range("S111").Value = 511
$ws.Range("T111").Value = 450
